$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the FX pair name, spot rate and volatility description for the
# new currency pair (USDEUR) as part of adding the CMS spread cap/floor.
$ws.Range("B2").Value = "USDEUR"
$ws.Range("B3").Value = 0.8384
$ws.Range("B4").Value = "USDEUR Vol 14Y coterm 3perc 30112020"

# Move the active selection, matching the author's final cursor position.
$ws.Range("F8").Select()
